$p = $ppt.ActivePresentation

# Slide 3 - Shape 2 ("Phan vai:" role assignments)
$s3 = $p.Slides.Item(3)
$sh1 = $s3.Shapes.Item(2)
$tr1 = $sh1.TextFrame.TextRange
    $tr1.Text = "Phân vai:"
    $tr1.InsertAfter("`r1. ") | Out-Null
    $tr1.InsertAfter("Nguyễn") | Out-Null
    $tr1.InsertAfter(" Minh ") | Out-Null
    $tr1.InsertAfter("Đức") | Out-Null
    $tr1.InsertAfter(": ") | Out-Null
    $tr1.InsertAfter("Đức") | Out-Null
    $tr1.InsertAfter(" ") | Out-Null
    $tr1.InsertAfter("lúc") | Out-Null
    $tr1.InsertAfter(" ") | Out-Null
    $tr1.InsertAfter("đã") | Out-Null
    $tr1.InsertAfter(" ") | Out-Null
    $tr1.InsertAfter("lớn") | Out-Null
    $tr1.InsertAfter(".") | Out-Null
    $tr1.InsertAfter("`r2. ") | Out-Null
    $tr1.InsertAfter("Trưởng") | Out-Null
    $tr1.InsertAfter(" ") | Out-Null
    $tr1.InsertAfter("Văn") | Out-Null
    $tr1.InsertAfter(" Khan: ") | Out-Null
    $tr1.InsertAfter("Đức") | Out-Null
    $tr1.InsertAfter(" ") | Out-Null
    $tr1.InsertAfter("lúc") | Out-Null
    $tr1.InsertAfter(" ") | Out-Null
    $tr1.InsertAfter("còn") | Out-Null
    $tr1.InsertAfter(" ") | Out-Null
    $tr1.InsertAfter("nhỏ") | Out-Null
    $tr1.InsertAfter(".") | Out-Null
    $tr1.InsertAfter("`r3. ") | Out-Null
    $tr1.InsertAfter("Nguyễn") | Out-Null
    $tr1.InsertAfter(" ") | Out-Null
    $tr1.InsertAfter("Ngọc") | Out-Null
    $tr1.InsertAfter(" ") | Out-Null
    $tr1.InsertAfter("Hảo") | Out-Null
    $tr1.InsertAfter(": ") | Out-Null
    $tr1.InsertAfter("`rCha ") | Out-Null
    $tr1.InsertAfter("của") | Out-Null
    $tr1.InsertAfter(" ") | Out-Null
    $tr1.InsertAfter("Đức") | Out-Null
    $tr1.InsertAfter(".") | Out-Null
    $tr1.InsertAfter("`r4. ") | Out-Null
    $tr1.InsertAfter("Nguyễn") | Out-Null
    $tr1.InsertAfter(" ") | Out-Null
    $tr1.InsertAfter("Tăng") | Out-Null
    $tr1.InsertAfter(" ") | Out-Null
    $tr1.InsertAfter("Tài") | Out-Null
    $tr1.InsertAfter(":") | Out-Null
    $tr1.InsertAfter("`rChú") | Out-Null
    $tr1.InsertAfter(" ") | Out-Null
    $tr1.InsertAfter("của") | Out-Null
    $tr1.InsertAfter(" ") | Out-Null
    $tr1.InsertAfter("Đức") | Out-Null
    $tr1.InsertAfter(".") | Out-Null
    $tr1.InsertAfter("`r5. ") | Out-Null
    $tr1.InsertAfter("Bùi") | Out-Null
    $tr1.InsertAfter(" ") | Out-Null
    $tr1.InsertAfter("Tiến") | Out-Null
    $tr1.InsertAfter(" ") | Out-Null
    $tr1.InsertAfter("Phát") | Out-Null
    $tr1.InsertAfter(": ") | Out-Null
    $tr1.InsertAfter("bạn") | Out-Null
    $tr1.InsertAfter(" ") | Out-Null
    $tr1.InsertAfter("học") | Out-Null
    $tr1.InsertAfter(" ") | Out-Null
    $tr1.InsertAfter("với") | Out-Null
    $tr1.InsertAfter(" ") | Out-Null
    $tr1.InsertAfter("Đức") | Out-Null
    $tr1.InsertAfter(" ") | Out-Null
    $tr1.InsertAfter("lúc") | Out-Null
    $tr1.InsertAfter(" ") | Out-Null
    $tr1.InsertAfter("lớn") | Out-Null
    $tr1.InsertAfter(".") | Out-Null
$sh1.Width = 330.0
$sh1.Height = 203.5687431574803

# Slide 3 - Shape 3 ("Dao cu + phan mem:" equipment list)
$sh2 = $s3.Shapes.Item(3)
$tr2 = $sh2.TextFrame.TextRange
    $tr2.Text = "Đạo cụ + phần mềm:"
    $tr2.InsertAfter("`r1. ") | Out-Null
    $tr2.InsertAfter("Máy") | Out-Null
    $tr2.InsertAfter(" Fujifilm ") | Out-Null
    $tr2.InsertAfter("Finepix") | Out-Null
    $tr2.InsertAfter(" S9400w") | Out-Null
    $tr2.InsertAfter("`r2. ") | Out-Null
    $tr2.InsertAfter("Hitfilm") | Out-Null
    $tr2.InsertAfter(" 4 Express") | Out-Null
    $tr2.InsertAfter("`r3. Audacity") | Out-Null
$sh2.Width = 264.0
$sh2.Height = 94.51409748818898

# Slide 7 - Shape 2 ("Nguoi long tieng:" voice actors)
$s7 = $p.Slides.Item(7)
$sh3 = $s7.Shapes.Item(2)
$tr3 = $sh3.TextFrame.TextRange
    $tr3.Text = "Người lồng tiếng:"
    $tr3.InsertAfter("`r1. ") | Out-Null
    $tr3.InsertAfter("Nguyễn") | Out-Null
    $tr3.InsertAfter(" Minh ") | Out-Null
    $tr3.InsertAfter("Đức") | Out-Null
    $tr3.InsertAfter(".") | Out-Null
    $tr3.InsertAfter("`r2. ") | Out-Null
    $tr3.InsertAfter("Nguyễn") | Out-Null
    $tr3.InsertAfter(" ") | Out-Null
    $tr3.InsertAfter("Tăng") | Out-Null
    $tr3.InsertAfter(" ") | Out-Null
    $tr3.InsertAfter("Tài") | Out-Null
    $tr3.InsertAfter(".") | Out-Null
    $tr3.InsertAfter("`r3. ") | Out-Null
    $tr3.InsertAfter("Bùi") | Out-Null
    $tr3.InsertAfter(" ") | Out-Null
    $tr3.InsertAfter("Tiến") | Out-Null
    $tr3.InsertAfter(" ") | Out-Null
    $tr3.InsertAfter("Phát") | Out-Null
    $tr3.InsertAfter(".") | Out-Null
    $tr3.InsertAfter("`r4. ") | Out-Null
    $tr3.InsertAfter("Nguyễn") | Out-Null
    $tr3.InsertAfter(" ") | Out-Null
    $tr3.InsertAfter("Ngọc") | Out-Null
    $tr3.InsertAfter(" ") | Out-Null
    $tr3.InsertAfter("Hảo") | Out-Null
    $tr3.InsertAfter(".") | Out-Null
    $tr3.InsertAfter("`r5. ") | Out-Null
    $tr3.InsertAfter("Trưởng") | Out-Null
    $tr3.InsertAfter(" ") | Out-Null
    $tr3.InsertAfter("Văn") | Out-Null
    $tr3.InsertAfter(" Khan") | Out-Null
$sh3.Width = 228.0
$sh3.Height = 138.135908511811

# Slide 7 - Shape 3 ("Nguoi bien tap:" editors)
$sh4 = $s7.Shapes.Item(3)
$tr4 = $sh4.TextFrame.TextRange
    $tr4.Text = "Người biên tập:"
    $tr4.InsertAfter("`r1. ") | Out-Null
    $tr4.InsertAfter("Nguyễn") | Out-Null
    $tr4.InsertAfter(" Minh ") | Out-Null
    $tr4.InsertAfter("Đức") | Out-Null
    $tr4.InsertAfter(".") | Out-Null
    $tr4.InsertAfter("`r2. ") | Out-Null
    $tr4.InsertAfter("Bùi") | Out-Null
    $tr4.InsertAfter(" ") | Out-Null
    $tr4.InsertAfter("Tiến") | Out-Null
    $tr4.InsertAfter(" ") | Out-Null
    $tr4.InsertAfter("Phát") | Out-Null
    $tr4.InsertAfter(".") | Out-Null
$sh4.Width = 324.0
$sh4.Height = 72.70315260629921

